$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eigen_Edelstahl")
$ws.Activate()

# Insert a new row at row 8 (pushes old row 8 -> row 9)
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the discount variable definition
$ws.Cells.Item(8, 1).Value = "Zahl"
$ws.Cells.Item(8, 2).Value = "Rabatt in %"
$ws.Cells.Item(8, 3).Value = "p_rabatt"

# Update the final price formula (now on row 9) to apply the discount
$ws.Cells.Item(9, 5).Value = "((L * P_Modell) + ((math.ceil(L/1.2)+1) * (P_Steher + P_Montageart)) + (Ecken * 150) + (L * P_Arbeit)) * (p_rabatt / 100)"

$ws.Range("D16").Select()
